$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the ID (column C) for row 24 (previously Michael Makris' ID)
$ws.Range("C24").ClearContents()

# Add new worker row 39: Michael Maksymciw, Number 1080, ID D4C6CF96
$ws.Range("A39").Value = "Michael Maksymciw"
$ws.Range("B39").Value = 1080
$ws.Range("C39").Value = "D4C6CF96"

# Update active selection as in the edited file
$ws.Range("I30").Select()
